# Auto-generated edit script to apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.506.90"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.562.69"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'211.67"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'46.01"
$ws.Range("E8").Value = "  +5.26%  "
$ws.Range("D9").Value = "'24.10"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'0.0591"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'0.0883"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "1.784.03"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "1.585.17"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'3.68"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "28.493.28"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'61.94"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").Value = "'227.06"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D21").Value = "'7.32"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  -6.40%  "
$ws.Range("D24").Value = "'9.11"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("E25").Value = "  +6.48%  "
$ws.Range("D26").Value = "'150.20"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("D27").Value = "'14.96"
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "'6.43"
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").Value = "'3.20"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "1.395.90"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'2.58"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "'0.536"
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").Value = "'5.53"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'62.80"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "1.697.23"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").Value = "'86.02"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  -4.32%  "
$ws.Range("E51").Value = "  -1.23%  "

Write-Output "Applied cryptos list update"
